$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"

$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "MCT-2A-M.T.R.M."
$ws.Range("E3").Value = "[-, -, -, 'MCT-3A-Elementos de máquinas']"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "MCT-2A-M.T.R.M."
$ws.Range("E4").Value = "[-, -, -, 'MCT-3A-Elementos de máquinas']"

$ws.Range("E6").Value = "[-, -, -, 'MCT-3A-Elementos de máquinas']"

$ws.Range("E7").Value = "[-, -, -, 'MCT-3A-Elementos de máquinas']"

$ws.Range("C8").Value = "-"
$ws.Range("F8").Value = "-"
